# Add a new data row (row 3) to the review sheet, mirroring the layout/
# formatting of the existing row 2, and add the accompanying hyperlink for
# the new email address in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone row 2's formatting onto row 3 first (keeps the same style indices
# used by the appid / email / recovery / time / review columns).
[void]$ws.Range("A2:F2").Copy()
[void]$ws.Range("A3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new row's values.
$ws.Range("A3").Value = "com.hamxa.shaynachim"
$ws.Range("B3").Value = "bitcoin guide"
$ws.Range("C3").Value = "ronenchen27@gmail.com"
$ws.Range("D3").Value = "danfogel100@gmail.com"
$ws.Range("E3").Value = "27/5/2019 15:59"
$ws.Range("F3").Value = "great beginners guide app. Makes a lot of sense in bitcoin"

# Add the mailto hyperlink on the new email cell (C3), matching the
# existing hyperlinks on C2/D2.
[void]$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:ronenchen27@gmail.com", [type]::Missing, [type]::Missing, "ronenchen27@gmail.com")

# Adding the hyperlink re-styles the cell with Excel's built-in "Hyperlink"
# look; restore the formatting that was copied from C2 so C3 keeps looking
# like the rest of the table.
$ws.Range("C3").Font.Name = "Calibri"
$ws.Range("C3").Font.Size = 11
$ws.Range("C3").Font.Color = 0
$ws.Range("C3").Font.Underline = -4142
$ws.Range("C3").HorizontalAlignment = -4108

# Match the saved selection state from the edit (B3 active cell).
[void]$ws.Range("B3").Select()
